$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 7) with an incomplete address (no Street Address Line 1/2)
# for the testing dataset.

$ws.Range("A7").Value = "Elementary School"
$ws.Range("B7").Value = "Nutcracker Academy"
$ws.Range("C7").Value = 1009876548
$ws.Range("D7").Value = "Ratty"
$ws.Range("E7").Value = "Teetherson"
$ws.Range("F7").Value = 17

# Date of Birth - stored as a real date value (unlike the text dates above it)
$ws.Range("G7").Value = 39872
$ws.Range("G7").NumberFormat = "mm-dd-yy"

# H7/I7 (Street Address Line 1/2) intentionally left blank - incomplete address

$ws.Range("J7").Value = "Hazelton"
$ws.Range("K7").Value = "Ontario"
$ws.Range("L7").Value = "N4U2L1"
$ws.Range("M7").Value = "HPV,"
$ws.Range("N7").Value = "HPV-9,"
$ws.Range("O7").Value = "Mar 12, 2014 - DTaP-IPV-Hib, Mar 12, 2014 - rota-unspecified, May 14, 2014 - Pneu-C-13, Jul 19, 2014 - DTaP-IPV-Hib, Sep 21, 2014 - MMR, Nov 25, 2014 - Men-C-C, Apr 17, 2015 - Var, Sep 13, 2015 - DTaP-IPV-Hib, May 5, 2024 - Tdap-IPV,"
$ws.Range("P7").Value = "NUTCRACKER ACADEMY-1009876547"
$ws.Range("Q7").Value = "HPV (HPV-9)"
$ws.Range("R7").Value = "[2014 MAR 12: DTaP-IPV-Hib, rota-unspecified] [2014 MAY 14: Pneu-C-13] [2014 JUL 19: DTaP-IPV-Hib] [2014 SEP 21: MMR] [2014 NOV 25: Men-C-C] [2015 APR 17: Var] [2015 SEP 13: DTaP-IPV-Hib] [2024 MAY 05: Tdap-IPV]"

# Auto-fit the Date of Birth column now that it holds a real formatted date
$ws.Range("G7").EntireColumn.AutoFit()

# Leave the selection where the edit happened
$ws.Range("H7").Select()
